$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the sigma_had value for the "O" (oxygen) row with the new current estimate.
$newValue = [double]"7.2585480000000004E+29"
$ws.Range("C2").Value = $newValue

# Move the active selection from F10 to F2 (formulas in D2/F2 recalc automatically).
$ws.Range("F2").Select()
